$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate metrics after trade #12 closed
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1399.66   # Current Capital
$wsSummary.Range("B4").Value = -0.35     # Total P&L $
$wsSummary.Range("B5").Value = -0.58     # Total P&L %
$wsSummary.Range("B6").Value = 12        # Total Trades
$wsSummary.Range("B8").Value = 9         # Losing Trades
$wsSummary.Range("B9").Value = 25        # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 5) metrics
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 99.66      # Capital
$wsStatus.Range("D5").Value = 12         # Trades
$wsStatus.Range("E5").Value = -0.35      # P&L $
$wsStatus.Range("F5").Value = -0.34      # P&L %
$wsStatus.Range("G5").Value = 25         # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly closed trade (#12) to both the "All Trades" and
# "MarketMaking" trade logs (row 13).
# ---------------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = 13

    $ws.Cells.Item($r, 1).Value = 12                 # Trade #

    # Date / Time are stored as plain text in this workbook (like the rows
    # above them) -- force a text number format first so Excel does not
    # reinterpret the literal as a date/time serial value.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = "2026-02-17"        # Date
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = "20:03:10"          # Time

    $ws.Cells.Item($r, 4).Value = "MarketMaking"      # Strategy
    $ws.Cells.Item($r, 5).Value = "DOWN"              # Side
    $ws.Cells.Item($r, 6).Value = 0.26                # Entry Price
    $ws.Cells.Item($r, 7).Value = 0.25                # Exit Price
    $ws.Cells.Item($r, 8).Value = "CLOSED"            # Status
    $ws.Cells.Item($r, 9).Value = -3.8462             # P&L %
    $ws.Cells.Item($r, 10).Value = -0.01              # P&L $
    $ws.Cells.Item($r, 11).Value = 99.66              # Capital After
    $ws.Cells.Item($r, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item($r, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item($r, 14).Value = 0.6                # Confidence
    $ws.Cells.Item($r, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($r, 16).Value = "early_exit"       # Exit Reason
    $ws.Cells.Item($r, 17).Value = 0.14               # Duration (min)
}
